$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B, C, D, E are treated as text so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00", "0.0000269") into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.846.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.303.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.00%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.89"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.299.24"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.94%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.53"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "630.35"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.831.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.12"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.840.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.58%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.297.65"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.906"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.87"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.61"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.98"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.38%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.93"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.65"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.18"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.07"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.48%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "544.82"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.105"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.795.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.49"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0739"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.58"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.42%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.129"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.20%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.23"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -12.75%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.73%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.80%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0416"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.129"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.96%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.12%  "
